# The document was re-saved through LibreOffice Writer with the
# document/paragraph language switched to English (USA). The visible
# text stays the same; what changes is:
#   - the run(s) carrying the sentence get merged into a single run
#     (LibreOffice's OOXML export coalesces runs / drops the
#     w:proofErr spell-check bookmarks Word had inserted), and
#   - that run (and, implicitly, the paragraph) is now tagged with
#     the English (USA) language instead of being language-less.
#
# Reproduce this with the Word object model: re-assert the paragraph
# text (which normalizes/merges the runs, since Word regenerates a
# single run for the replacement) and then switch the language of the
# whole body to English (US). This also picks up the "cs" (complex
# script) font that LibreOffice stamps onto the run.

$d = $word.ActiveDocument

$text = "FORMATIERUNG: Calibri/Carlito, 12pt und 1,5 Zeilenabstand (ist voreingestellt)"

# Re-assert the paragraph's text via Find/Replace across the whole
# body. Word/this engine collapses the matched (multi-run) span into a
# single run carrying the first run's formatting, which also drops the
# w:proofErr spell-check markers that surrounded "Carlito".
$null = $d.Content.Find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, $text, 2)

# Switch the document body's language to English (USA), matching the
# "Sprachen ... auf Englisch (USA)" language-settings change.
$d.Content.LanguageID = "en-US"

# LibreOffice also stamps the complex-script (w:cs) font explicitly as
# Calibri on the run; mirror that via the Bi (bidi/complex-script) font
# name.
$d.Content.Font.NameBi = "Calibri"
